$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be auto-parsed as numbers
foreach ($r in @(4,5,6,7,8,9,10,11,13,14,15,16,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,40,41,43,45,46,48,49,50,51)) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply cell updates per the diff
$ws.Range("D2").Value = "26.121.34"
$ws.Range("D3").Value = "1.670.29"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "216.59"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "0.5217"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.2696"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "0.06393"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "21.85"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "0.07446"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "1.690.88"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "4.521"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "0.5825"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "0.000008519"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "64.17"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "25.907.18"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "4.934"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "10.80"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "189.39"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").Value = "6.198"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("D24").Value = "144.70"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "0.1244"
$ws.Range("E25").Value = "  +5.32%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "7.607"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "15.75"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "0.06572"
$ws.Range("E28").Value = "  +13.14%  "
$ws.Range("D29").Value = "1.343"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "1.316"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").Value = "3.587"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "3.528"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "1.666"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").Value = "1.019"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "0.6173"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").Value = "2.368"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Value = "2.697"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").Value = "6.257"
$ws.Range("E38").Value = "  +6.35%  "
$ws.Range("D39").Value = "1.094.43"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "0.01597"
$ws.Range("D41").Value = "0.8692"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "100.91"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").Value = "1.816.80"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("D46").Value = "56.50"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "8.134"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").Value = "0.05237"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "0.4277"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "5.996"
$ws.Range("E51").Value = "  +2.22%  "

Write-Host "Applied cryptos update"